# MENT-169: Create New Question Categories
#
# Adds a new "Oficial de Pediatria" entry under the existing
# "CLINICAL_AREA" / "Área Clinica" category, inserted right after the
# last existing CLINICAL_AREA row (old row 97, "Oficial de ATS") and
# before the "MONITORING_AND_EVALUATION" section that used to start at
# row 98. Inserting the row pushes the three MONITORING_AND_EVALUATION
# rows down by one (old rows 98-100 -> new rows 99-101) and Excel keeps
# their B/C/D contents intact automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 98; everything from row 98 down (including
# the three MONITORING_AND_EVALUATION rows) shifts down by one.
$ws.Rows("98:98").Insert()

# Fill in the new row with the new "Oficial de Pediatria" category entry.
$ws.Range("A98").Value = 96
$ws.Range("B98").Value = "CLINICAL_AREA"
$ws.Range("C98").Value = "Área Clinica"
$ws.Range("D98").Value = "Oficial de Pediatria"

# The "Nr" column (A) holds literal sequence numbers, not a formula, so
# renumber the rows that shifted down.
$ws.Range("A99").Value = 97
$ws.Range("A100").Value = 98
$ws.Range("A101").Value = 99

# Match the final saved selection/viewport from the authored workbook.
$ws.Range("B20").Select()
